{"js": "// Fix duplicate ETests section.\n//\n// The template has two unrelated placeholders that happen to share the\n// name \"ETests\":\n//   - an outer section wrapper: {#ETests} ... {/ETests}\n//     (small 6pt runs, wraps \"Antibiotikaempfindlichkeitstestung ...\"\n//     and the sentence ending in \"{/ETests}{Interpretation}\")\n//   - an inner per-item loop: {#ETests}{Antibiotic} ... {ValidFromYear}{/ETests}\n//     (11pt runs)\n//\n// Because both use the identical tag name \"ETests\", the outer section is\n// ambiguous/duplicated against the inner loop. The fix renames only the\n// OUTER wrapper tag to \"HasETests\", leaving the inner loop's \"ETests\"\n// markers untouched.\n//\n// We find the two paragraphs that hold the outer wrapper's open/close\n// tags by their (unique) leading text, then replace just the \"ETests\"\n// run within each of those paragraphs.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst openParaStart = \"{#ETests}Antibiotikaempfindlichkeitstestung\";\nconst closeParaStart = \"{/ETests}{Interpretation}\";\n\nconst targetParagraphs = paragraphs.items.filter((p) =>\n  p.text.indexOf(openParaStart) === 0 || p.text.indexOf(closeParaStart) === 0\n);\n\nconst searchResults = targetParagraphs.map((p) =>\n  p.search(\"ETests\", { matchCase: true, matchWholeWord: true })\n);\nawait context.sync();\n\nfor (const results of searchResults) {\n  for (const item of results.items) {\n    item.insertText(\"HasETests\", Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Fix duplicate ETests section.\n#\n# The template has two unrelated placeholders that happen to share the\n# name \"ETests\":\n#   - an outer section wrapper: {#ETests} ... {/ETests}\n#     (wraps \"Antibiotikaempfindlichkeitstestung ...\" and the sentence\n#     ending in \"{/ETests}{Interpretation}\")\n#   - an inner per-item loop: {#ETests}{Antibiotic} ... {ValidFromYear}{/ETests}\n#\n# Because both use the identical tag name \"ETests\", the outer section is\n# ambiguous/duplicated against the inner loop. The fix renames only the\n# OUTER wrapper tag to \"HasETests\", leaving the inner loop's \"ETests\"\n# markers untouched.\n#\n# We locate the two paragraphs that hold the outer wrapper's open/close\n# tags by their (unique) leading text, then replace just the \"ETests\"\n# run within each of those paragraphs.\n\n$d = $word.ActiveDocument\n\n$openParaStart = \"{#ETests}Antibiotikaempfindlichkeitstestung\"\n$closeParaStart = \"{/ETests}{Interpretation}\"\n\nforeach ($p in $d.Paragraphs) {\n    $paraText = $p.Range.Text\n    if ($paraText.StartsWith($openParaStart) -or $paraText.StartsWith($closeParaStart)) {\n        $scoped = $p.Range.Duplicate\n        $scoped.Find.ClearFormatting()\n        $scoped.Find.Text = \"ETests\"\n        $scoped.Find.MatchCase = $true\n        $scoped.Find.MatchWholeWord = $true\n        $scoped.Find.Forward = $true\n        $scoped.Find.Wrap = 0\n        if ($scoped.Find.Execute()) {\n            $scoped.Text = \"HasETests\"\n        }\n    }\n}\n"}
